$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.724.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.730.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4936'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2624'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06219'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.727.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '15.93'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6118'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.507'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9985'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.522.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9983'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007215'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.948.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.488'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.574'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.104'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.773'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.387'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.938'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07986'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.675'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04486'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9978'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.609'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6247'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9424'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.043'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.422'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01513'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.581'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3861'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.950'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1160'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05382'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.869'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.36%  '
